# Update the "取得日時" (fetched timestamp) column on the "ランサーズ" sheet
# for all existing data rows (2-7) to the new run timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-14 06:40:02"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
